$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 74.566666
$ws.Range("I5").Value = 68
$ws.Range("J5").Value = 133.66667
$ws.Range("K5").Value = 68
$ws.Range("L5").Value = 133.66667
$ws.Range("M5").Value = 47
$ws.Range("N5").Value = -363.66667
$ws.Range("H69").Value = 4994
$ws.Range("I69").Value = 4994
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 14982
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -14108
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 4994
$ws.Range("I72").Value = 4994
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 44946
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -40578
$ws.Range("N72").ClearContents()
$ws.Range("H98").Value = 5124.2324
$ws.Range("I98").Value = 4835.5312
$ws.Range("K98").Value = 4835.5312
$ws.Range("M98").Value = -3337.5312
$ws.Range("H100").Value = 1963.1111
$ws.Range("I100").Value = 1337.375
$ws.Range("J100").Value = 6969
$ws.Range("K100").Value = 1337.375
$ws.Range("L100").Value = 6969
$ws.Range("M100").Value = -796.375
$ws.Range("N100").Value = -8051
$ws.Range("H103").Value = 554.8889
$ws.Range("J103").Value = 596.1539
$ws.Range("L103").Value = 1788.4617
$ws.Range("N103").Value = -2960.4617
$ws.Range("H122").Value = 5124.2324
$ws.Range("I122").Value = 4835.5312
$ws.Range("K122").Value = 14506.5936
$ws.Range("M122").Value = -12056.5936
$ws.Range("H125").Value = 41667476
$ws.Range("J125").Value = 800
$ws.Range("L125").Value = 7200
$ws.Range("N125").Value = -12120
$ws.Range("H137").Value = 3927.7
$ws.Range("I137").Value = 6069.25
$ws.Range("J137").Value = 2500
$ws.Range("K137").Value = 18207.75
$ws.Range("L137").Value = 7500
$ws.Range("M137").Value = -15657.75
$ws.Range("N137").Value = -12600

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1787597.4
$ws.Range("I32").Value = 1813433.5
$ws.Range("K32").Value = 1813433.5
$ws.Range("M32").Value = -1813146.5
$ws.Range("H61").Value = 6064.04
$ws.Range("I61").Value = 3437.303
$ws.Range("K61").Value = 3437.303
$ws.Range("M61").Value = -3225.303
$ws.Range("H63").Value = 2731.6667
$ws.Range("I63").Value = 2597.5
$ws.Range("K63").Value = 2597.5
$ws.Range("M63").Value = -1911.5
$ws.Range("H66").Value = 2731.6667
$ws.Range("I66").Value = 2597.5
$ws.Range("K66").Value = 12987.5
$ws.Range("M66").Value = -9555.5
$ws.Range("H136").Value = 6064.04
$ws.Range("I136").Value = 3437.303
$ws.Range("K136").Value = 10311.909
$ws.Range("M136").Value = -7761.909

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 6894.4
$ws.Range("I22").Value = 8570.083000000001
$ws.Range("K22").Value = 8570.083000000001
$ws.Range("M22").Value = -8397.083000000001
$ws.Range("H134").Value = 7131.485
$ws.Range("I134").Value = 3649.6667
$ws.Range("K134").Value = 10949.0001
$ws.Range("M134").Value = -8414.000100000001

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5086.3076
$ws.Range("I16").Value = 2806.4
$ws.Range("J16").Value = 6511.25
$ws.Range("K16").Value = 2806.4
$ws.Range("L16").Value = 6511.25
$ws.Range("M16").Value = -2519.4
$ws.Range("N16").Value = -7085.25
$ws.Range("H58").Value = 10422208
$ws.Range("I58").Value = 29413410
$ws.Range("K58").Value = 29413410
$ws.Range("M58").Value = -29413207
$ws.Range("H105").Value = 7144997
$ws.Range("I105").Value = 10205182
$ws.Range("J105").Value = 4566
$ws.Range("K105").Value = 10205182
$ws.Range("L105").Value = 4566
$ws.Range("M105").Value = -10203435
$ws.Range("N105").Value = -8060
$ws.Range("H113").Value = 5086.3076
$ws.Range("I113").Value = 2806.4
$ws.Range("J113").Value = 6511.25
$ws.Range("K113").Value = 2806.4
$ws.Range("L113").Value = 6511.25
$ws.Range("M113").Value = -636.4000000000001
$ws.Range("N113").Value = -10851.25
$ws.Range("H132").Value = 4439.4575
$ws.Range("I132").Value = 2822.625
$ws.Range("K132").Value = 8467.875
$ws.Range("M132").Value = -5937.875
$ws.Range("H136").Value = 10422208
$ws.Range("I136").Value = 29413410
$ws.Range("K136").Value = 88240230
$ws.Range("M136").Value = -88237680

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 28.8
$ws.Range("J38").Value = 29.333334
$ws.Range("L38").Value = 88.00000199999999
$ws.Range("N38").Value = -782.000002
$ws.Range("H48").Value = 8332.333000000001
$ws.Range("J48").Value = 8332.333000000001
$ws.Range("L48").Value = 24996.999
$ws.Range("N48").Value = -25496.999
$ws.Range("H60").Value = 2575
$ws.Range("I60").Value = 2575
$ws.Range("K60").Value = 7725
$ws.Range("M60").Value = -7474
$ws.Range("H122").Value = 2022517.6
$ws.Range("I122").Value = 4042683
$ws.Range("J122").Value = 2352.1428
$ws.Range("K122").Value = 36384147
$ws.Range("L122").Value = 21169.2852
$ws.Range("M122").Value = -36381697
$ws.Range("N122").Value = -26069.2852
$ws.Range("H131").Value = 1913.1052
$ws.Range("I131").Value = 890.1
$ws.Range("J131").Value = 3049.7778
$ws.Range("K131").Value = 2670.3
$ws.Range("L131").Value = 9149.3334
$ws.Range("M131").Value = 2369.7
$ws.Range("N131").Value = -19229.3334

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4123.364
$ws.Range("I7").Value = 3883.6667
$ws.Range("K7").Value = 3883.6667
$ws.Range("M7").Value = -3771.6667
$ws.Range("H126").Value = 4123.364
$ws.Range("I126").Value = 3883.6667
$ws.Range("K126").Value = 11651.0001
$ws.Range("M126").Value = -9181.000100000001
$ws.Range("H136").Value = 14606.286
$ws.Range("I136").Value = 5744
$ws.Range("K136").Value = 17232
$ws.Range("M136").Value = -14682

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 18188392
$ws.Range("I81").Value = 1366.7778
$ws.Range("K81").Value = 2733.5556
$ws.Range("M81").Value = -1672.5556
$ws.Range("H84").Value = 18188392
$ws.Range("I84").Value = 1366.7778
$ws.Range("K84").Value = 13667.778
$ws.Range("M84").Value = -8363.778
$ws.Range("H132").Value = 31289872
$ws.Range("I132").Value = 50012148
$ws.Range("J132").Value = 86073.164
$ws.Range("K132").Value = 150036444
$ws.Range("L132").Value = 258219.492
$ws.Range("M132").Value = -150033914
$ws.Range("N132").Value = -263279.492
